# "Actualizar" run: appends one new availability-check pass (14 services)
# to the bottom of the log and nudges the timestamp recorded for the
# previous pass (tiny float re-quantization, same as the source automation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-stamp the last existing pass (rows 198-211) with the
#        slightly adjusted timestamp seen in the committed workbook.
$prevTimestamp = 44231.96609321759
for ($r = 198; $r -le 211; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $prevTimestamp
}

# --- 2. Append the new pass: 14 services, same order as every prior
#        block, all stamped with the new run's timestamp.
$newTimestamp = 44231.98718089004

$services = @(
    @{Row=212; Name="Odoo";               DisplayUrl="https://www.dataintelligence-group.com/"; Url="https://www.dataintelligence-group.com/"},
    @{Row=213; Name="Blackbox";            DisplayUrl="https://serviciodashboard.azurewebsites.net/"; Url="https://serviciodashboard.azurewebsites.net/"},
    @{Row=214; Name="PowerBI";             DisplayUrl="https://powerbi.microsoft.com/es-es/"; Url="https://powerbi.microsoft.com/es-es/"},
    @{Row=215; Name="Dropbox";             DisplayUrl="https://www.dropbox.com/"; Url="https://www.dropbox.com/"},
    @{Row=216; Name="Odoo";                DisplayUrl="https://dataintelligence.store/"; Url="https://dataintelligence.store/"},
    @{Row=217; Name="GEE";                 DisplayUrl="https://app-data-i.users.earthengine.app/"; Url="https://app-data-i.users.earthengine.app/"},
    @{Row=218; Name="UtilidadesOdoo";      DisplayUrl="https://odooutil.azurewebsites.net/"; Url="https://odooutil.azurewebsites.net/"},
    @{Row=219; Name="Filtros Dashboard";   DisplayUrl="https://filtradordashboard.azurewebsites.net/"; Url="https://filtradordashboard.azurewebsites.net/"},
    @{Row=220; Name="MapStore";            DisplayUrl="https://ide.dataintelligence-group.com/mapstore/#/"; Url="https://ide.dataintelligence-group.com/mapstore/"; SubAddress="/"},
    @{Row=221; Name="GeoServer";           DisplayUrl="https://ide.dataintelligence-group.com/geoserver/web/?0"; Url="https://ide.dataintelligence-group.com/geoserver/web/?0"},
    @{Row=222; Name="Tomcat";              DisplayUrl="https://ide.dataintelligence-group.com/"; Url="https://ide.dataintelligence-group.com/"},
    @{Row=223; Name="Shiny";               DisplayUrl="https://rpubs.com/dataintelligence/"; Url="https://rpubs.com/dataintelligence/"},
    @{Row=224; Name="Github";              DisplayUrl="https://github.com/Sud-Austral/"; Url="https://github.com/Sud-Austral/"},
    @{Row=225; Name="EZ Exporter";         DisplayUrl="https://ezexporter.highviewapps.com/exports/export-profile/"; Url="https://ezexporter.highviewapps.com/exports/export-profile/"}
)

foreach ($svc in $services) {
    $r = $svc.Row

    $ws.Cells.Item($r, 1).Value = $svc.Name

    # Set the visible text first so Hyperlinks.Add (which can split off a
    # '#fragment' into SubAddress) doesn't change what the cell displays.
    $linkCell = $ws.Cells.Item($r, 2)
    $linkCell.Value = $svc.DisplayUrl
    if ($svc.ContainsKey("SubAddress")) {
        $ws.Hyperlinks.Add($linkCell, $svc.Url, $svc.SubAddress)
    } else {
        $ws.Hyperlinks.Add($linkCell, $svc.Url)
    }
    $linkCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $dateCell = $ws.Cells.Item($r, 4)
    $dateCell.Value2 = $newTimestamp
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
